$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 51
$ws1.Range("F6").Value = 9375
$ws1.Range("F7").Value = 841
$ws1.Range("F10").Value = 1127
$ws1.Range("F12").Value = 92
$ws1.Range("F14").Value = 259
$ws1.Range("F15").Value = 409
$ws1.Range("F16").Value = 87
$ws1.Range("F18").Value = 1265

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 51
$ws4.Range("F7").Value = 9375
$ws4.Range("F8").Value = 841
$ws4.Range("F11").Value = 1127
$ws4.Range("F13").Value = 92
$ws4.Range("F15").Value = 259
$ws4.Range("F16").Value = 409
$ws4.Range("F17").Value = 87
$ws4.Range("F19").Value = 1265
